$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns per latest data refresh.
# Values are written as text to match the existing inline-string cell format
# (NumberFormat is temporarily set to Text so Excel does not auto-convert
# numeric-looking strings into numbers/dates), then the cell style is reset
# back to Normal so no extra formatting is introduced.
$updates = @(
    @('D2', '28.933.16'),
    @('D3', '1.901.44'),
    @('E3', '  -4.75%  '),
    @('E4', '  +0.28%  '),
    @('D5', '323.46'),
    @('E5', '  -1.19%  '),
    @('D6', '1.003'),
    @('E6', '  +0.15%  '),
    @('D7', '0.4592'),
    @('E7', '  -2.00%  '),
    @('D8', '0.3803'),
    @('E8', '  -3.61%  '),
    @('D9', '45.52'),
    @('E9', '  -2.40%  '),
    @('D10', '0.07689'),
    @('E10', '  -4.34%  '),
    @('D11', '0.9787'),
    @('E11', '  -2.19%  '),
    @('D12', '22.01'),
    @('E12', '  -3.77%  '),
    @('D13', '1.893.62'),
    @('E13', '  -4.91%  '),
    @('D14', '6.944'),
    @('E14', '  -4.12%  '),
    @('D15', '5.654'),
    @('E15', '  -3.55%  '),
    @('D16', '0.07043'),
    @('E16', '  -1.23%  '),
    @('D17', '1.005'),
    @('E17', '  +0.13%  '),
    @('D18', '83.60'),
    @('E18', '  -6.10%  '),
    @('D19', '0.000009496'),
    @('E19', '  -5.43%  '),
    @('D20', '16.66'),
    @('E20', '  -4.41%  '),
    @('E21', '  +0.21%  '),
    @('D22', '28.941.73'),
    @('E22', '  -2.15%  '),
    @('D23', '5.307'),
    @('E23', '  -4.69%  '),
    @('D24', '10.86'),
    @('E24', '  -3.44%  '),
    @('D25', '2.127.06'),
    @('E25', '  -4.78%  '),
    @('D26', '2.093'),
    @('E26', '  -0.43%  '),
    @('D27', '157.10'),
    @('E27', '  -0.64%  '),
    @('D28', '19.06'),
    @('E28', '  -3.34%  '),
    @('D29', '5.576'),
    @('E29', '  -6.58%  '),
    @('D30', '117.26'),
    @('E30', '  -2.61%  '),
    @('D31', '1.845'),
    @('E31', '  -5.12%  '),
    @('D32', '0.09274'),
    @('E32', '  -1.93%  '),
    @('D33', '0.8588'),
    @('E33', '  -6.39%  '),
    @('D34', '5.059'),
    @('E34', '  -4.27%  '),
    @('D35', '1.243'),
    @('E35', '  -8.84%  '),
    @('D36', '3.029'),
    @('E36', '  -5.90%  '),
    @('D37', '0.05694'),
    @('E37', '  -2.43%  '),
    @('D38', '1.140'),
    @('E38', '  -3.20%  '),
    @('E39', '  +0.17%  '),
    @('E40', '  -4.19%  '),
    @('E41', '  -5.74%  '),
    @('D42', '0.5496'),
    @('E42', '  -4.80%  '),
    @('E43', '  -3.85%  '),
    @('D44', '9.225'),
    @('E44', '  -6.73%  '),
    @('D45', '2.754'),
    @('E45', '  -1.79%  '),
    @('D46', '0.5188'),
    @('E46', '  -4.00%  '),
    @('D48', '2.091'),
    @('E48', '  -5.72%  '),
    @('E49', '  -2.48%  '),
    @('D50', '111.27'),
    @('E50', '  -2.48%  '),
    @('D51', '1.775'),
    @('E51', '  -5.05%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
    $rng.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates."
